$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240, shifting rows 240:250 down to 241:251.
$ws.Rows("240:240").Insert()

# Populate the newly inserted row 240 with the new data record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,T are constant across the whole data block
# (same market / product), so copy them from the row immediately above (239).
$ws.Range("A240").Value = $ws.Range("A239").Value2
$ws.Range("B240").Value = $ws.Range("B239").Value2
$ws.Range("C240").Value = $ws.Range("C239").Value2
$ws.Range("D240").Value = 44826
$ws.Range("E240").Value = $ws.Range("E239").Value2
$ws.Range("F240").Value = $ws.Range("F239").Value2
$ws.Range("G240").Value = $ws.Range("G239").Value2
$ws.Range("H240").Value = $ws.Range("H239").Value2
$ws.Range("I240").Value = $ws.Range("I239").Value2
$ws.Range("J240").Value = $ws.Range("J239").Value2
$ws.Range("K240").Value = $ws.Range("K239").Value2
$ws.Range("L240").Value = "Primera"
$ws.Range("M240").Value = 200
$ws.Range("N240").Value = 10000
$ws.Range("O240").Value = 11000
$ws.Range("P240").Value = 10500
$ws.Range("Q240").Value = $ws.Range("Q239").Value2
$ws.Range("R240").Value = "Brasil"
$ws.Range("S240").Value = 2625
$ws.Range("T240").Value = $ws.Range("T239").Value2
